$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 26234.75
$ws.Range("J17").Value = 26234.75
$ws.Range("L17").Value = 78704.25
$ws.Range("N17").Value = -79040.25
$ws.Range("H58").Value = 18690.443
$ws.Range("I58").Value = 286.66666
$ws.Range("J58").Value = 21875.71
$ws.Range("K58").Value = 859.9999799999999
$ws.Range("L58").Value = 65627.13
$ws.Range("M58").Value = -709.9999799999999
$ws.Range("N58").Value = -65927.13
$ws.Range("H69").Value = 4464
$ws.Range("I69").Value = 3813
$ws.Range("J69").Value = 5115
$ws.Range("K69").Value = 11439
$ws.Range("L69").Value = 15345
$ws.Range("M69").Value = -10565
$ws.Range("N69").Value = -17093
$ws.Range("H72").Value = 4464
$ws.Range("I72").Value = 3813
$ws.Range("J72").Value = 5115
$ws.Range("K72").Value = 34317
$ws.Range("L72").Value = 46035
$ws.Range("M72").Value = -29949
$ws.Range("N72").Value = -54771
$ws.Range("H96").Value = 624.82355
$ws.Range("I96").Value = 390.44446
$ws.Range("J96").Value = 888.5
$ws.Range("K96").Value = 1171.33338
$ws.Range("L96").Value = 2665.5
$ws.Range("M96").Value = 201.66662
$ws.Range("N96").Value = -5411.5
$ws.Range("H103").Value = 427.8125
$ws.Range("I103").Value = 700
$ws.Range("J103").Value = 388.92856
$ws.Range("K103").Value = 2100
$ws.Range("L103").Value = 1166.78568
$ws.Range("M103").Value = -1514
$ws.Range("N103").Value = -2338.78568
$ws.Range("H106").Value = 676.55554
$ws.Range("I106").Value = 398.42856
$ws.Range("J106").Value = 1650
$ws.Range("K106").Value = 398.42856
$ws.Range("L106").Value = 1650
$ws.Range("M106").Value = 232.57144
$ws.Range("N106").Value = -2912
$ws.Range("H115").Value = 1208.75
$ws.Range("I115").Value = 611.6667
$ws.Range("K115").Value = 1835.0001
$ws.Range("M115").Value = -268.0001
$ws.Range("H132").Value = 6064367
$ws.Range("I132").Value = 9528077
$ws.Range("J132").Value = 2874.5833
$ws.Range("K132").Value = 28584231
$ws.Range("L132").Value = 8623.749899999999
$ws.Range("M132").Value = -28581701
$ws.Range("N132").Value = -13683.7499
$ws.Range("H137").Value = 3139.282
$ws.Range("I137").Value = 3104.0312
$ws.Range("J137").Value = 3300.4285
$ws.Range("K137").Value = 9312.0936
$ws.Range("L137").Value = 9901.2855
$ws.Range("M137").Value = -6762.0936
$ws.Range("N137").Value = -15001.2855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1026.0625
$ws.Range("I74").Value = 730.53845
$ws.Range("J74").Value = 2306.6667
$ws.Range("K74").Value = 730.53845
$ws.Range("L74").Value = 2306.6667
$ws.Range("M74").Value = 143.46155
$ws.Range("N74").Value = -4054.6667
$ws.Range("H77").Value = 1026.0625
$ws.Range("I77").Value = 730.53845
$ws.Range("J77").Value = 2306.6667
$ws.Range("K77").Value = 3652.69225
$ws.Range("L77").Value = 11533.3335
$ws.Range("M77").Value = 715.3077499999999
$ws.Range("N77").Value = -20269.3335
$ws.Range("H132").Value = 2122.2827
$ws.Range("I132").Value = 1699.625
$ws.Range("J132").Value = 3088.3572
$ws.Range("K132").Value = 5098.875
$ws.Range("L132").Value = 9265.071599999999
$ws.Range("M132").Value = -2568.875
$ws.Range("N132").Value = -14325.0716

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1021.5333
$ws.Range("I94").Value = 915.3
$ws.Range("J94").Value = 1234
$ws.Range("K94").Value = 915.3
$ws.Range("L94").Value = 1234
$ws.Range("M94").Value = -464.3
$ws.Range("N94").Value = -2136
$ws.Range("H134").Value = 2125
$ws.Range("I134").Value = 1369.4445
$ws.Range("J134").Value = 5525
$ws.Range("K134").Value = 4108.333500000001
$ws.Range("L134").Value = 16575
$ws.Range("M134").Value = -1573.333500000001
$ws.Range("N134").Value = -21645

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H120").Value = 21750
$ws.Range("J120").Value = 21750
$ws.Range("L120").Value = 21750
$ws.Range("N120").Value = -29008
$ws.Range("H132").Value = 2057.3096
$ws.Range("I132").Value = 1617.5518
$ws.Range("J132").Value = 3038.3076
$ws.Range("K132").Value = 4852.6554
$ws.Range("L132").Value = 9114.9228
$ws.Range("M132").Value = -2322.6554
$ws.Range("N132").Value = -14174.9228

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 3823.75
$ws.Range("J63").Value = 3941.4285
$ws.Range("L63").Value = 11824.2855
$ws.Range("N63").Value = -13322.2855
$ws.Range("H66").Value = 3823.75
$ws.Range("J66").Value = 3941.4285
$ws.Range("L66").Value = 35472.8565
$ws.Range("N66").Value = -42960.8565
$ws.Range("H87").Value = 6804.4
$ws.Range("I87").Value = 2307.8572
$ws.Range("J87").Value = 12527.272
$ws.Range("K87").Value = 6923.571599999999
$ws.Range("L87").Value = 37581.81600000001
$ws.Range("M87").Value = -5675.571599999999
$ws.Range("N87").Value = -40077.81600000001
$ws.Range("H90").Value = 6804.4
$ws.Range("I90").Value = 2307.8572
$ws.Range("J90").Value = 12527.272
$ws.Range("K90").Value = 20770.7148
$ws.Range("L90").Value = 112745.448
$ws.Range("M90").Value = -14530.7148
$ws.Range("N90").Value = -125225.448
$ws.Range("H120").Value = 18757.875
$ws.Range("J120").Value = 18838.834
$ws.Range("L120").Value = 56516.50199999999
$ws.Range("N120").Value = -66192.50199999999
$ws.Range("H122").Value = 1227.8
$ws.Range("I122").Value = 465.8
$ws.Range("J122").Value = 1989.8
$ws.Range("K122").Value = 4192.2
$ws.Range("L122").Value = 17908.2
$ws.Range("M122").Value = -1742.2
$ws.Range("N122").Value = -22808.2
$ws.Range("H124").Value = 12808.462
$ws.Range("J124").Value = 17786.666
$ws.Range("L124").Value = 53359.99800000001
$ws.Range("N124").Value = -63179.99800000001
$ws.Range("H140").Value = 20836180
$ws.Range("I140").Value = 33334106
$ws.Range("K140").Value = 100002318
$ws.Range("M140").Value = -99997138

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3845.4517
$ws.Range("I132").Value = 4032.0527
$ws.Range("J132").Value = 3550
$ws.Range("K132").Value = 12096.1581
$ws.Range("L132").Value = 10650
$ws.Range("M132").Value = -9566.158100000001
$ws.Range("N132").Value = -15710

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2109.0908
$ws.Range("I46").Value = 400
$ws.Range("J46").Value = 2488.889
$ws.Range("K46").Value = 400
$ws.Range("L46").Value = 2488.889
$ws.Range("M46").Value = -212
$ws.Range("N46").Value = -2864.889
$ws.Range("H132").Value = 2204.4102
$ws.Range("I132").Value = 1332.4
$ws.Range("K132").Value = 3997.2
$ws.Range("M132").Value = -1467.2
$ws.Range("H136").Value = 2523.7446
$ws.Range("I136").Value = 1385.1578
$ws.Range("J136").Value = 7331.1113
$ws.Range("K136").Value = 4155.4734
$ws.Range("L136").Value = 21993.3339
$ws.Range("M136").Value = -1605.4734
$ws.Range("N136").Value = -27093.3339

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 10423.667
$ws.Range("I132").Value = 2063.75
$ws.Range("J132").Value = 33413.438
$ws.Range("K132").Value = 6191.25
$ws.Range("L132").Value = 100240.314
$ws.Range("M132").Value = -3661.25
$ws.Range("N132").Value = -105300.314
